$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "14/05/2024"
$ws.Range("C2").Value = "15/05/2024"
$ws.Range("D2").Value = "Diana Gómez"
$ws.Range("E2").Value = "Moto Yamaha FZ"
$ws.Range("F2").Value = "Florencia (Caquetá)"
$ws.Range("G2").Value = "Terminal"
$ws.Range("H2").Value = "Gilberto Gómez"
$ws.Range("I2").Value = 220000
$ws.Range("J2").Value = 33000

$ws.Range("B3").Value = "10/05/2024"
$ws.Range("C3").Value = "15/05/2024"
$ws.Range("D3").Value = "Mario Gómez"
$ws.Range("E3").Value = "Moto Yamaha FZ"
$ws.Range("F3").Value = "Montería (Córdoba)"
$ws.Range("G3").Value = "Terminal"
$ws.Range("H3").Value = "Gilberto Gómez"
$ws.Range("I3").Value = 370000
$ws.Range("J3").Value = 55500

$ws.Range("B4").Value = "28/04/2024"
$ws.Range("C4").Value = "29/04/2024"
$ws.Range("D4").Value = "Diana Gómez"
$ws.Range("E4").Value = "Carro Toyota Land Cruiser Prado"
$ws.Range("F4").Value = "Santa Marta (Magdalena)"
$ws.Range("G4").Value = "Centro Comercial"
$ws.Range("H4").Value = "Mario Gómez"
$ws.Range("I4").Value = 380000
$ws.Range("J4").Value = 57000

$ws.Range("B5").Value = "11/04/2024"
$ws.Range("C5").Value = "23/04/2024"
$ws.Range("D5").Value = "Diana Gómez"
$ws.Range("E5").Value = "Moto Yamaha XTZ"
$ws.Range("F5").Value = "Arauca (Arauca)"
$ws.Range("G5").Value = "Terminal"
$ws.Range("H5").Value = "Diana Caicedo"
$ws.Range("I5").Value = 285000
$ws.Range("J5").Value = 42750

$ws.Range("B6").Value = "10/04/2024"
$ws.Range("C6").Value = "15/04/2024"
$ws.Range("D6").Value = "Mario Gómez"
$ws.Range("E6").Value = "Moto Yamaha FZ"
$ws.Range("F6").Value = "Florencia (Caquetá)"
$ws.Range("G6").Value = "Terminal"
$ws.Range("H6").Value = "Diana Caicedo"
$ws.Range("I6").Value = 120000
$ws.Range("J6").Value = 18000

$ws.Range("B7").Value = "26/03/2024"
$ws.Range("C7").Value = "29/03/2024"
$ws.Range("D7").Value = "Diana Gómez"
$ws.Range("E7").Value = "Carro Toyota Land Cruiser Prado"
$ws.Range("F7").Value = "Santa Marta (Magdalena)"
$ws.Range("G7").Value = "Terminal"
$ws.Range("H7").Value = "Mario Gómez"
$ws.Range("I7").Value = 300000
$ws.Range("J7").Value = 45000

$ws.Range("B8").Value = "20/03/2024"
$ws.Range("C8").Value = "21/03/2024"
$ws.Range("D8").Value = "Mario Gómez"
$ws.Range("E8").Value = "Moto Yamaha FZ"
$ws.Range("F8").Value = "Montería (Córdoba)"
$ws.Range("G8").Value = "Aeropuerto"
$ws.Range("H8").Value = "Diana Caicedo"
$ws.Range("I8").Value = 90000
$ws.Range("J8").Value = 13500

$ws.Range("B9").Value = "24/02/2024"
$ws.Range("C9").Value = "25/02/2024"
$ws.Range("D9").Value = "Mario Gómez"
$ws.Range("E9").Value = "Carro Chevrolet Spark"
$ws.Range("F9").Value = "Florencia (Caquetá)"
$ws.Range("G9").Value = "Aeropuerto"
$ws.Range("H9").Value = "Gilberto Gómez"
$ws.Range("I9").Value = 400000
$ws.Range("J9").Value = 60000

$ws.Range("B10").Value = "24/02/2024"
$ws.Range("C10").Value = "25/02/2024"
$ws.Range("D10").Value = "Mario Gómez"
$ws.Range("E10").Value = "Moto Honda Eco Deluxe"
$ws.Range("F10").Value = "Florencia (Caquetá)"
$ws.Range("G10").Value = "Terminal"
$ws.Range("H10").Value = "Diana Caicedo"
$ws.Range("I10").Value = 215000
$ws.Range("J10").Value = 32250

$ws.Range("B11").Value = "17/02/2024"
$ws.Range("C11").Value = "19/02/2024"
$ws.Range("D11").Value = "Mario Gómez"
$ws.Range("E11").Value = "Moto Honda XR"
$ws.Range("F11").Value = "Cali (Valle del Cauca)"
$ws.Range("G11").Value = "Aeropuerto"
$ws.Range("H11").Value = "Gilberto Gómez"
$ws.Range("I11").Value = 390000
$ws.Range("J11").Value = 58500

$ws.Range("B12").Value = "8/02/2024"
$ws.Range("C12").Value = "9/02/2024"
$ws.Range("D12").Value = "Diana Gómez"
$ws.Range("E12").Value = "Moto Yamaha FZ"
$ws.Range("F12").Value = "Florencia (Caquetá)"
$ws.Range("G12").Value = "Aeropuerto"
$ws.Range("H12").Value = "Gilberto Gómez"
$ws.Range("I12").Value = 380000
$ws.Range("J12").Value = 57000

$ws.Range("B13").Value = "1/02/2024"
$ws.Range("C13").Value = "9/02/2024"
$ws.Range("D13").Value = "Mario Gómez"
$ws.Range("E13").Value = "Moto Yamaha FZ"
$ws.Range("F13").Value = "Montería (Córdoba)"
$ws.Range("G13").Value = "Terminal"
$ws.Range("H13").Value = "Diana Caicedo"
$ws.Range("I13").Value = 725000
$ws.Range("J13").Value = 108750
